# Auto-generated edit script
# Applies row-content rotations/swaps among rows 15-20 and 24-26
# of worksheet "Artfynd", matching the target OOXML diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A15").Value = "130938723"
$ws.Range("AC15").Value = ""
$ws.Range("AF15").Value = ""
$ws.Range("AF15").NumberFormat = "General"
$ws.Range("AM15").Value = ""
$ws.Range("AO15").Value = "Picea abies"
$ws.Range("B15").Value = "83223"
$ws.Range("E15").Value = "6440"
$ws.Range("F15").Value = "Vitgrynig nållav"
$ws.Range("G15").Value = "Chaenotheca subroscida"
$ws.Range("H15").Value = "(Eitner) Zahlbr."
$ws.Range("J15").Value = ""
$ws.Range("J15").NumberFormat = "General"
$ws.Range("L15").Value = ""
$ws.Range("M15").Value = ""
$ws.Range("Q15").Value = "476289"
$ws.Range("R15").Value = "7033519"
$ws.Range("A16").Value = "130938729"
$ws.Range("AC16").Value = "Ringhack, äldre, på gran."
$ws.Range("AF16").Value = ""
$ws.Range("AM16").Value = "Trädstam på levande träd"
$ws.Range("AO16").Value = "Stem on living tree # Picea abies"
$ws.Range("B16").Value = "57884"
$ws.Range("E16").Value = "100109"
$ws.Range("F16").Value = "Tretåig hackspett"
$ws.Range("G16").Value = "Picoides tridactylus"
$ws.Range("H16").Value = "(Linnaeus, 1758)"
$ws.Range("J16").Value = ""
$ws.Range("L16").Value = ""
$ws.Range("L16").NumberFormat = "General"
$ws.Range("M16").Value = "äldre spår"
$ws.Range("Q16").Value = "476618"
$ws.Range("R16").Value = "7033500"
$ws.Range("A17").Value = "130938740"
$ws.Range("AJ17").Value = ""
$ws.Range("AK17").Value = ""
$ws.Range("AO17").Value = ""
$ws.Range("B17").Value = "78255"
$ws.Range("E17").Value = "228579"
$ws.Range("F17").Value = "Liten svartspik"
$ws.Range("G17").Value = "Chaenothecopsis nana"
$ws.Range("H17").Value = "Tibell"
$ws.Range("K17").Value = ""
$ws.Range("K17").NumberFormat = "General"
$ws.Range("Q17").Value = "476532"
$ws.Range("R17").Value = "7033586"
$ws.Range("A18").Value = "130938743"
$ws.Range("AJ18").Value = "gran"
$ws.Range("AK18").Value = "Picea abies"
$ws.Range("AO18").Value = "Picea abies"
$ws.Range("B18").Value = "91828"
$ws.Range("E18").Value = "5432"
$ws.Range("F18").Value = "Granticka"
$ws.Range("G18").Value = "Porodaedalea chrysoloma s.lat."
$ws.Range("H18").Value = ""
$ws.Range("H18").NumberFormat = "General"
$ws.Range("K18").Value = "teleomorf"
$ws.Range("Q18").Value = "476555"
$ws.Range("R18").Value = "7033581"
$ws.Range("A19").Value = "130938734"
$ws.Range("AC19").Value = "Ringhack, äldre, på gran."
$ws.Range("AF19").Value = ""
$ws.Range("AM19").Value = "Trädstam på levande träd"
$ws.Range("AO19").Value = "Stem on living tree # Picea abies"
$ws.Range("B19").Value = "57884"
$ws.Range("E19").Value = "100109"
$ws.Range("F19").Value = "Tretåig hackspett"
$ws.Range("G19").Value = "Picoides tridactylus"
$ws.Range("H19").Value = "(Linnaeus, 1758)"
$ws.Range("J19").Value = ""
$ws.Range("L19").Value = ""
$ws.Range("L19").NumberFormat = "General"
$ws.Range("M19").Value = "äldre spår"
$ws.Range("Q19").Value = "476457"
$ws.Range("R19").Value = "7033634"
$ws.Range("A20").Value = "130938752"
$ws.Range("AC20").Value = ""
$ws.Range("AF20").Value = ""
$ws.Range("AF20").NumberFormat = "General"
$ws.Range("AM20").Value = ""
$ws.Range("AO20").Value = "Picea abies"
$ws.Range("B20").Value = "79243"
$ws.Range("E20").Value = "6425"
$ws.Range("F20").Value = "Garnlav"
$ws.Range("G20").Value = "Alectoria sarmentosa"
$ws.Range("H20").Value = "(Ach.) Ach."
$ws.Range("J20").Value = ""
$ws.Range("J20").NumberFormat = "General"
$ws.Range("L20").Value = ""
$ws.Range("M20").Value = ""
$ws.Range("Q20").Value = "476286"
$ws.Range("R20").Value = "7033527"
$ws.Range("A24").Value = "130938751"
$ws.Range("AC24").Value = ""
$ws.Range("AF24").Value = ""
$ws.Range("AF24").NumberFormat = "General"
$ws.Range("B24").Value = "79243"
$ws.Range("E24").Value = "6425"
$ws.Range("F24").Value = "Garnlav"
$ws.Range("G24").Value = "Alectoria sarmentosa"
$ws.Range("H24").Value = "(Ach.) Ach."
$ws.Range("J24").Value = ""
$ws.Range("J24").NumberFormat = "General"
$ws.Range("L24").Value = ""
$ws.Range("M24").Value = ""
$ws.Range("Q24").Value = "476394"
$ws.Range("R24").Value = "7033617"
$ws.Range("A25").Value = "130938746"
$ws.Range("Q25").Value = "476419"
$ws.Range("R25").Value = "7033605"
$ws.Range("A26").Value = "130938731"
$ws.Range("AC26").Value = "Ringhack, äldre, på gran."
$ws.Range("AF26").Value = ""
$ws.Range("B26").Value = "57884"
$ws.Range("E26").Value = "100109"
$ws.Range("F26").Value = "Tretåig hackspett"
$ws.Range("G26").Value = "Picoides tridactylus"
$ws.Range("H26").Value = "(Linnaeus, 1758)"
$ws.Range("J26").Value = ""
$ws.Range("L26").Value = ""
$ws.Range("L26").NumberFormat = "General"
$ws.Range("M26").Value = "äldre spår"
$ws.Range("Q26").Value = "476549"
$ws.Range("R26").Value = "7033604"
